$d = $word.ActiveDocument

# Locate the paragraph that still holds the Word field code
# ( { m:'Mona_Lisa.jpg'.asImage().fit(100, 400, false) } ) built out of
# fldChar/instrText runs, and rewrite it as plain literal-text runs
# ( w:t ) using "{" / "}" instead of the field delimiters, exactly as
# produced by TokenIteratorFieldRewriterSplit.

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the paragraph containing the field code"
}

$color = '<w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr>'

$runs = ''
$runs += '<w:r><w:t>{</w:t></w:r>'
$runs += '<w:r><w:t>m</w:t></w:r>'
$runs += '<w:r><w:t>:</w:t></w:r>'
$runs += '<w:r>' + $color + "<w:t>'</w:t></w:r>"
$runs += '<w:r>' + $color + '<w:t>Mona_Lisa.jpg</w:t></w:r>'
$runs += '<w:r>' + $color + "<w:t>'.asImage()</w:t></w:r>"
$runs += '<w:r>' + $color + '<w:t>.fit(</w:t></w:r>'
$runs += '<w:r>' + $color + '<w:t>1</w:t></w:r>'
$runs += '<w:r>' + $color + '<w:t>0</w:t></w:r>'
$runs += '<w:r>' + $color + '<w:t xml:space="preserve">0, </w:t></w:r>'
$runs += '<w:r>' + $color + '<w:t>400</w:t></w:r>'
$runs += '<w:r>' + $color + '<w:t>, false</w:t></w:r>'
$runs += '<w:r>' + $color + '<w:t>)</w:t></w:r>'
$runs += '<w:r><w:t xml:space="preserve">}</w:t></w:r>'

$paraAttrs = 'w14:paraId="2F8A187F" w14:textId="180EC619" w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F"'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
    '<w:body><w:p ' + $paraAttrs + '>' + $runs + '</w:p></w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$start = $target.Range.Start
$end = $target.Range.End
$r = $d.Range($start, $end - 1)
$r.InsertXML($xml)
